$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing 을지대학교(성남) row (row 20/21 merged block) ---
$ws.Range("E20").Value = "특성화고특별전형(학종)"
$ws.Range("H20").Value = "(스마트의료정부학부)의료공학전공"
$ws.Range("R20").Value = 2.91
$ws.Range("S20").Value = 3.89
$ws.Range("T20").Value = 3.36

# --- 2. Add new 서울신학대학교 row (row 24/25 merged block) ---
$ws.Range("B24").Value = "서울신학대학교"
$ws.Range("E24").Value = "(교과)`n특성화고교졸업자전형(정원 외)"
$ws.Range("H24").Value = "IT융합소프트웨어학과"
$ws.Range("K24").Value = "6(2?)"
$ws.Range("L24").Value = "해당없음"
$ws.Range("N24").Value = "국어, 수학, 영어, 사회, 과학 *상위 3과목*`n특성화고교, 예술고, 체육고 및 학력인정교(2년제, 방송`n고, 대안교육특성화학교), 직업교육 및 대안교육 위탁과`n정은 전체 교과목 반영"
$ws.Range("Q24").Value = "857.58`n(3.85)"

# Row 25 needs a taller custom height to fit the long text in row 24/25's merged block
$ws.Rows(25).RowHeight = 117

# --- 3. Column J (10) width change ---
$ws.Columns(10).ColumnWidth = 19.75

# --- 4. Sheet view: scroll position + selection ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("S24:S25").Select
